# Add data for 2021-10-30 (updates the "through October 2x" running-total
# column plus a handful of individual neighborhood/month cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the header label for the "current month" column.
$ws.Name = "Through 2021-10-22"
$ws.Range("B1").Value = "October 2021 (through October 22)"

# Updated counts (existing cells whose value increased).
$ws.Range("AZ2").Value = 5
$ws.Range("AZ3").Value = 2
$ws.Range("B4").Value = 10
$ws.Range("B9").Value = 6
$ws.Range("AZ13").Value = 3
$ws.Range("V24").Value = 2
$ws.Range("B27").Value = 3
$ws.Range("B66").Value = 4
$ws.Range("L66").Value = 3

# New counts (cells that were previously empty).
$ws.Range("V8").Value = 1
$ws.Range("AZ20").Value = 1
$ws.Range("AF21").Value = 1
$ws.Range("AP33").Value = 2
$ws.Range("B59").Value = 1
$ws.Range("AP60").Value = 1
$ws.Range("AP67").Value = 1
$ws.Range("AP87").Value = 1
$ws.Range("L93").Value = 1
